# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with refreshed figures, as produced by the scheduled GitHub Actions scrape.
#
# All cells in D/E are plain text in the source workbook (prices such as
# "64.269.79" use '.' as both thousands separators, which Excel cannot
# parse as a real number anyway; percentages keep their padding spaces).
# Whenever a new price happens to look like a genuine number (e.g. "590.92"),
# a direct .Value assignment would make Excel auto-convert it to a numeric
# cell (losing trailing zeros / switching to scientific notation) and that
# also reformats the cell. To keep every D/E cell a plain string with its
# original (default) formatting untouched, those values are written via a
# text formula and then flattened to a static value with Copy +
# PasteSpecial(values only) instead of being assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.344.61'
$ws.Range('E2').Value = '  +8.70%  '
$ws.Range('D3').Value = '3.158.40'
$ws.Range('E3').Value = '  +6.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Formula = '="590.92"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +4.42%  '
$ws.Range('D6').Formula = '="148.16"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +8.21%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.152.17'
$ws.Range('E8').Value = '  +6.04%  '
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('D10').Formula = '="0.159"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +20.27%  '
$ws.Range('E11').Value = '  +10.00%  '
$ws.Range('D12').Formula = '="0.473"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +5.29%  '
$ws.Range('D13').Formula = '="0.0000256"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +12.20%  '
$ws.Range('D14').Formula = '="35.98"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +7.29%  '
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '3.681.16'
$ws.Range('E16').Value = '  +6.14%  '
$ws.Range('D17').Value = '64.232.45'
$ws.Range('E17').Value = '  +8.53%  '
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').Value = '3.154.94'
$ws.Range('E19').Value = '  +6.15%  '
$ws.Range('D20').Formula = '="477.40"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +9.78%  '
$ws.Range('D21').Formula = '="14.29"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +4.94%  '
$ws.Range('D22').Formula = '="0.736"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('E23').Value = '  +9.52%  '
$ws.Range('D24').Formula = '="13.44"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('D25').Formula = '="82.75"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +13.32%  '
$ws.Range('E28').Value = '  +6.66%  '
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Formula = '="6.92"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +12.36%  '
$ws.Range('D32').Formula = '="27.38"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +6.80%  '
$ws.Range('E33').Value = '  +6.54%  '
$ws.Range('D34').Value = '0.0₃0892'
$ws.Range('E34').Value = '  +16.61%  '
$ws.Range('D35').Formula = '="2.44"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +18.42%  '
$ws.Range('E36').Value = '  +7.65%  '
$ws.Range('D37').Formula = '="3.43"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +23.21%  '
$ws.Range('D38').Formula = '="6.18"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +5.23%  '
$ws.Range('D39').Formula = '="50.97"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +5.17%  '
$ws.Range('D40').Formula = '="452.94"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +14.84%  '
$ws.Range('D41').Formula = '="8.80"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('E42').Value = '  +6.90%  '
$ws.Range('D43').Value = '2.941.80'
$ws.Range('E43').Value = '  +8.43%  '
$ws.Range('D44').Formula = '="0.285"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +14.38%  '
$ws.Range('E45').Value = '  +6.70%  '
$ws.Range('D46').Formula = '="2.24"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +13.61%  '
$ws.Range('D47').Formula = '="35.68"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +3.98%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Formula = '="123.45"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('D51').Formula = '="25.10"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +8.50%  '

$excel.CutCopyMode = $false
